# The "startup" sheet (xl/worksheets/sheet1.xml) contains a small lookup
# table of query metadata. The "caseDetailQuery" column (column C) used a
# hard-coded case id (NCATS-COP01CCB010072) in its Cypher query, which is
# no longer needed/correct, so the whole column is removed, shifting the
# dbExcel/WebExcel columns (D, E) left into C, D.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire "caseDetailQuery" column (column C), shifting the
# remaining columns to the left - same as right-clicking the column C
# header in Excel and choosing "Delete".
$ws.Columns.Item(3).Delete()

# Select the new column C (mirrors the resulting selection left behind by
# Excel after a column delete).
$ws.Columns.Item(3).Select()
